$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates reflecting refreshed crypto price/volume data
$ws.Range("D2").Value = "64.373.37"
$ws.Range("E2").Value = "  -3.06%  "
$ws.Range("D3").Value = "3.181.63"
$ws.Range("E3").Value = "  -7.99%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.54"
$ws.Range("E5").Value = "  -3.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.27"
$ws.Range("E6").Value = "  -3.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.609"
$ws.Range("E7").Value = "  -0.85%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "3.181.63"
$ws.Range("E10").Value = "  -6.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.64"
$ws.Range("E12").Value = "  -4.85%  "
$ws.Range("D13").Value = "3.731.75"
$ws.Range("E13").Value = "  -8.02%  "
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.45"
$ws.Range("E15").Value = "  -9.34%  "
$ws.Range("D16").Value = "64.366.72"
$ws.Range("E16").Value = "  -2.91%  "
$ws.Range("E17").Value = "  -5.11%  "
$ws.Range("D18").Value = "3.181.77"
$ws.Range("E18").Value = "  -7.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.75"
$ws.Range("E19").Value = "  -3.74%  "
$ws.Range("E20").Value = "  -5.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "353.63"
$ws.Range("E21").Value = "  -5.42%  "
$ws.Range("E22").Value = "  -5.45%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  -5.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000120"
$ws.Range("E25").Value = "  -4.22%  "
$ws.Range("E26").Value = "  -5.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.56"
$ws.Range("E27").Value = "  -3.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.177"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.61"
$ws.Range("E30").Value = "  -4.62%  "
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("E32").Value = "  -4.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.09"
$ws.Range("E33").Value = "  -6.89%  "
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.65"
$ws.Range("E34").Value = "  -5.63%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.20"
$ws.Range("E35").Value = "  -5.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.44"
$ws.Range("E36").Value = "  -6.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "155.00"
$ws.Range("E37").Value = "  -3.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.821"
$ws.Range("E38").Value = "  -7.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.97"
$ws.Range("E39").Value = "  -8.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.56"
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("E41").Value = "  -5.69%  "
$ws.Range("D42").Value = "2.630.14"
$ws.Range("E42").Value = "  -4.93%  "
$ws.Range("E43").Value = "  -7.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.03"
$ws.Range("E44").Value = "  -6.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "39.61"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0658"
$ws.Range("E46").Value = "  -5.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.87"
$ws.Range("E47").Value = "  -5.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "323.63"
$ws.Range("E48").Value = "  -4.71%  "
$ws.Range("E49").Value = "  -7.40%  "
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  -0.01%  "
